$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.699.30'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.536.46'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.19'
$ws.Range("E5").Value = '  -2.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.74'
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E7").Value = '  -1.09%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.21'
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.34'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.931.26'
$ws.Range("E14").Value = '  -1.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.83'
$ws.Range("E15").Value = '  +4.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.544.93'
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.812'
$ws.Range("E17").Value = '  -3.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.674.70'
$ws.Range("E18").Value = '  -1.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.24'
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.38'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.10'
$ws.Range("E23").Value = '  -3.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.90'
$ws.Range("E24").Value = '  -2.55%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.05'
$ws.Range("E27").Value = '  -4.47%  '
$ws.Range("E28").Value = '  -5.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.23'
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.81'
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("E33").Value = '  +12.52%  '
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.04'
$ws.Range("E36").Value = '  -5.49%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.41'
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.17'
$ws.Range("E38").Value = '  -7.01%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.32'
$ws.Range("E41").Value = '  +8.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.34'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0299'
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.969.28'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '80.69'
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.727.38'
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("E51").Value = '  +9.10%  '
